$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 868369.7
$ws.Range("I6").Value = 1157648.5
$ws.Range("J6").Value = 533.3333
$ws.Range("K6").Value = 3472945.5
$ws.Range("L6").Value = 1599.9999
$ws.Range("M6").Value = -3472833.5
$ws.Range("N6").Value = -1823.9999

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 48.125
$ws.Range("I8").Value = 48.125
$ws.Range("K8").Value = 144.375
$ws.Range("M8").Value = -5.375

# ALC row 34
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 6918.5
$ws.Range("I34").Value = 6918.5
$ws.Range("K34").Value = 6918.5
$ws.Range("M34").Value = -6715.5

# ALC row 36
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 6918.5
$ws.Range("I36").Value = 6918.5
$ws.Range("K36").Value = 6918.5
$ws.Range("M36").Value = -6203.5

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 15880779
$ws.Range("I76").Value = 10280.714
$ws.Range("K76").Value = 10280.714
$ws.Range("M76").Value = -9965.714

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 15880779
$ws.Range("I79").Value = 10280.714
$ws.Range("K79").Value = 10280.714
$ws.Range("M79").Value = -9188.714

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3820.5618
$ws.Range("I138").Value = 1709.2858
$ws.Range("J138").Value = 5188.9814
$ws.Range("K138").Value = 5127.857400000001
$ws.Range("L138").Value = 15566.9442
$ws.Range("M138").Value = 12.14259999999922
$ws.Range("N138").Value = -25846.9442

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31280406
$ws.Range("I32").Value = 58839436
$ws.Range("K32").Value = 58839436
$ws.Range("M32").Value = -58839149

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2549.375
$ws.Range("I63").Value = 2549.375
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2549.375
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1863.375
$ws.Range("N63").ClearContents()

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2549.375
$ws.Range("I66").Value = 2549.375
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12746.875
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9314.875
$ws.Range("N66").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2825.9644
$ws.Range("I74").Value = 3519.5625
$ws.Range("J74").Value = 1901.1666
$ws.Range("K74").Value = 3519.5625
$ws.Range("L74").Value = 1901.1666
$ws.Range("M74").Value = -2645.5625
$ws.Range("N74").Value = -3649.1666

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2825.9644
$ws.Range("I77").Value = 3519.5625
$ws.Range("J77").Value = 1901.1666
$ws.Range("K77").Value = 17597.8125
$ws.Range("L77").Value = 9505.833000000001
$ws.Range("M77").Value = -13229.8125
$ws.Range("N77").Value = -18241.833

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 32321
$ws.Range("I122").Value = 32321
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 96963
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -94513
$ws.Range("N122").ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 27781866
$ws.Range("I132").Value = 43479324
$ws.Range("J132").Value = 9441.385
$ws.Range("K132").Value = 130437972
$ws.Range("L132").Value = 28324.155
$ws.Range("M132").Value = -130435442
$ws.Range("N132").Value = -33384.155

# BSM row 5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2675
$ws.Range("I5").Value = 173.33333
$ws.Range("J5").Value = 3508.889
$ws.Range("K5").Value = 173.33333
$ws.Range("L5").Value = 3508.889
$ws.Range("M5").Value = -60.33332999999999
$ws.Range("N5").Value = -3734.889

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 55556492
$ws.Range("I107").Value = 71429416
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 71429416
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = -71427496
$ws.Range("N107").Value = -5090

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2025781.8
$ws.Range("I134").Value = 4727.6177
$ws.Range("K134").Value = 14182.8531
$ws.Range("M134").Value = -11647.8531

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7579265
$ws.Range("I132").Value = 970.6923
$ws.Range("J132").Value = 18525690
$ws.Range("K132").Value = 2912.0769
$ws.Range("L132").Value = 55577070
$ws.Range("M132").Value = -382.0769
$ws.Range("N132").Value = -55582130

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 808.03
$ws.Range("I131").Value = 527.5
$ws.Range("J131").Value = 819.71875
$ws.Range("K131").Value = 1582.5
$ws.Range("L131").Value = 2459.15625
$ws.Range("M131").Value = 3457.5
$ws.Range("N131").Value = -12539.15625

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 811.8570999999999
$ws.Range("I97").Value = 822.1667
$ws.Range("K97").Value = 822.1667
$ws.Range("M97").Value = -326.1667

# GSM row 125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# LTW row 18
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6331113.5
$ws.Range("I22").Value = 6331113.5
$ws.Range("K22").Value = 6331113.5
$ws.Range("M22").Value = -6330818.5

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 6331113.5
$ws.Range("I27").Value = 6331113.5
$ws.Range("K27").Value = 6331113.5
$ws.Range("M27").Value = -6331006.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 13889318
$ws.Range("I46").Value = 13889318
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 13889318
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -13889130
$ws.Range("N46").ClearContents()

# LTW row 47
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 11000
$ws.Range("J47").Value = 11000
$ws.Range("L47").Value = 11000
$ws.Range("N47").Value = -11980

# LTW row 52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 11000
$ws.Range("J52").Value = 11000
$ws.Range("L52").Value = 11000
$ws.Range("N52").Value = -11466

# LTW row 62
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 37624.5
$ws.Range("J62").Value = 37624.5
$ws.Range("L62").Value = 37624.5
$ws.Range("N62").Value = -38872.5

# LTW row 65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 37624.5
$ws.Range("J65").Value = 37624.5
$ws.Range("L65").Value = 112873.5
$ws.Range("N65").Value = -119113.5

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2722

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -4496

# LTW row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 39000
$ws.Range("J108").Value = 39000
$ws.Range("L108").Value = 39000
$ws.Range("N108").Value = -46680

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 14677.647
$ws.Range("I122").Value = 20413.637
$ws.Range("J122").Value = 4161.6665
$ws.Range("K122").Value = 61240.91099999999
$ws.Range("L122").Value = 12484.9995
$ws.Range("M122").Value = -58790.91099999999
$ws.Range("N122").Value = -17384.9995

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17262266
$ws.Range("I132").Value = 41710584
$ws.Range("J132").Value = 4629.9414
$ws.Range("K132").Value = 125131752
$ws.Range("L132").Value = 13889.8242
$ws.Range("M132").Value = -125129222
$ws.Range("N132").Value = -18949.8242
